$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALERTS")
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = "2026-01-28"
$ws.Cells.Item(3,2).Value = "14:41:10"
$ws.Cells.Item(3,3).Value = "14:00"
$ws.Cells.Item(3,4).Value = "Bathroom"
$ws.Cells.Item(3,5).Value = "MINIMAL"
$ws.Cells.Item(3,6).Value = "MINIMAL ALERT: Bathroom occupied, no motion > 20s."
$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "2026-01-28"
$ws.Cells.Item(4,2).Value = "14:41:30"
$ws.Cells.Item(4,3).Value = "14:00"
$ws.Cells.Item(4,4).Value = "Bathroom"
$ws.Cells.Item(4,5).Value = "MODERATE"
$ws.Cells.Item(4,6).Value = "MODERATE ALERT: Bathroom occupied, no motion > 40s."

$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-01-28"
$ws.Cells.Item(14,2).Value = "14:40:35"
$ws.Cells.Item(14,3).Value = "14:00"
$ws.Cells.Item(14,4).Value = "Bathroom"
$ws.Cells.Item(14,5).Value = "No Motion"
$ws.Cells.Item(14,6).Value = "Inactive"
$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-01-28"
$ws.Cells.Item(15,2).Value = "14:40:37"
$ws.Cells.Item(15,3).Value = "14:00"
$ws.Cells.Item(15,4).Value = "Bathroom"
$ws.Cells.Item(15,5).Value = "No Motion"
$ws.Cells.Item(15,6).Value = "Inactive"
$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-01-28"
$ws.Cells.Item(16,2).Value = "14:40:42"
$ws.Cells.Item(16,3).Value = "14:00"
$ws.Cells.Item(16,4).Value = "Bathroom"
$ws.Cells.Item(16,5).Value = "No Motion"
$ws.Cells.Item(16,6).Value = "Inactive"
$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = "2026-01-28"
$ws.Cells.Item(17,2).Value = "14:40:47"
$ws.Cells.Item(17,3).Value = "14:00"
$ws.Cells.Item(17,4).Value = "Bathroom"
$ws.Cells.Item(17,5).Value = "No Motion"
$ws.Cells.Item(17,6).Value = "Inactive"
$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = "2026-01-28"
$ws.Cells.Item(18,2).Value = "14:40:52"
$ws.Cells.Item(18,3).Value = "14:00"
$ws.Cells.Item(18,4).Value = "Bathroom"
$ws.Cells.Item(18,5).Value = "No Motion"
$ws.Cells.Item(18,6).Value = "Inactive"
$ws.Cells.Item(19,1).NumberFormat = "@"
$ws.Cells.Item(19,1).Value = "2026-01-28"
$ws.Cells.Item(19,2).Value = "14:40:57"
$ws.Cells.Item(19,3).Value = "14:00"
$ws.Cells.Item(19,4).Value = "Bathroom"
$ws.Cells.Item(19,5).Value = "No Motion"
$ws.Cells.Item(19,6).Value = "Inactive"
$ws.Cells.Item(20,1).NumberFormat = "@"
$ws.Cells.Item(20,1).Value = "2026-01-28"
$ws.Cells.Item(20,2).Value = "14:41:02"
$ws.Cells.Item(20,3).Value = "14:00"
$ws.Cells.Item(20,4).Value = "Bathroom"
$ws.Cells.Item(20,5).Value = "No Motion"
$ws.Cells.Item(20,6).Value = "Inactive"
$ws.Cells.Item(21,1).NumberFormat = "@"
$ws.Cells.Item(21,1).Value = "2026-01-28"
$ws.Cells.Item(21,2).Value = "14:41:07"
$ws.Cells.Item(21,3).Value = "14:00"
$ws.Cells.Item(21,4).Value = "Bathroom"
$ws.Cells.Item(21,5).Value = "No Motion"
$ws.Cells.Item(21,6).Value = "Inactive"
$ws.Cells.Item(22,1).NumberFormat = "@"
$ws.Cells.Item(22,1).Value = "2026-01-28"
$ws.Cells.Item(22,2).Value = "14:41:12"
$ws.Cells.Item(22,3).Value = "14:00"
$ws.Cells.Item(22,4).Value = "Bathroom"
$ws.Cells.Item(22,5).Value = "No Motion"
$ws.Cells.Item(22,6).Value = "Inactive"
$ws.Cells.Item(23,1).NumberFormat = "@"
$ws.Cells.Item(23,1).Value = "2026-01-28"
$ws.Cells.Item(23,2).Value = "14:41:17"
$ws.Cells.Item(23,3).Value = "14:00"
$ws.Cells.Item(23,4).Value = "Bathroom"
$ws.Cells.Item(23,5).Value = "No Motion"
$ws.Cells.Item(23,6).Value = "Inactive"
$ws.Cells.Item(24,1).NumberFormat = "@"
$ws.Cells.Item(24,1).Value = "2026-01-28"
$ws.Cells.Item(24,2).Value = "14:41:22"
$ws.Cells.Item(24,3).Value = "14:00"
$ws.Cells.Item(24,4).Value = "Bathroom"
$ws.Cells.Item(24,5).Value = "No Motion"
$ws.Cells.Item(24,6).Value = "Inactive"
$ws.Cells.Item(25,1).NumberFormat = "@"
$ws.Cells.Item(25,1).Value = "2026-01-28"
$ws.Cells.Item(25,2).Value = "14:41:27"
$ws.Cells.Item(25,3).Value = "14:00"
$ws.Cells.Item(25,4).Value = "Bathroom"
$ws.Cells.Item(25,5).Value = "No Motion"
$ws.Cells.Item(25,6).Value = "Inactive"
$ws.Cells.Item(26,1).NumberFormat = "@"
$ws.Cells.Item(26,1).Value = "2026-01-28"
$ws.Cells.Item(26,2).Value = "14:41:32"
$ws.Cells.Item(26,3).Value = "14:00"
$ws.Cells.Item(26,4).Value = "Bathroom"
$ws.Cells.Item(26,5).Value = "No Motion"
$ws.Cells.Item(26,6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-01-28"
$ws.Cells.Item(14,2).Value = "14:40:35"
$ws.Cells.Item(14,3).Value = "14:00"
$ws.Cells.Item(14,4).Value = "Bathroom"
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "88.8%"
$ws.Cells.Item(14,6).Value = "Active"
$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-01-28"
$ws.Cells.Item(15,2).Value = "14:40:39"
$ws.Cells.Item(15,3).Value = "14:00"
$ws.Cells.Item(15,4).Value = "Bathroom"
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = "87.8%"
$ws.Cells.Item(15,6).Value = "Active"
$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-01-28"
$ws.Cells.Item(16,2).Value = "14:40:43"
$ws.Cells.Item(16,3).Value = "14:00"
$ws.Cells.Item(16,4).Value = "Bathroom"
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = "88.8%"
$ws.Cells.Item(16,6).Value = "Active"
$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = "2026-01-28"
$ws.Cells.Item(17,2).Value = "14:40:51"
$ws.Cells.Item(17,3).Value = "14:00"
$ws.Cells.Item(17,4).Value = "Bathroom"
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = "88.8%"
$ws.Cells.Item(17,6).Value = "Active"
$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = "2026-01-28"
$ws.Cells.Item(18,2).Value = "14:40:55"
$ws.Cells.Item(18,3).Value = "14:00"
$ws.Cells.Item(18,4).Value = "Bathroom"
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = "88.8%"
$ws.Cells.Item(18,6).Value = "Active"
$ws.Cells.Item(19,1).NumberFormat = "@"
$ws.Cells.Item(19,1).Value = "2026-01-28"
$ws.Cells.Item(19,2).Value = "14:40:59"
$ws.Cells.Item(19,3).Value = "14:00"
$ws.Cells.Item(19,4).Value = "Bathroom"
$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value = "87.8%"
$ws.Cells.Item(19,6).Value = "Active"
$ws.Cells.Item(20,1).NumberFormat = "@"
$ws.Cells.Item(20,1).Value = "2026-01-28"
$ws.Cells.Item(20,2).Value = "14:41:03"
$ws.Cells.Item(20,3).Value = "14:00"
$ws.Cells.Item(20,4).Value = "Bathroom"
$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value = "88.7%"
$ws.Cells.Item(20,6).Value = "Active"
$ws.Cells.Item(21,1).NumberFormat = "@"
$ws.Cells.Item(21,1).Value = "2026-01-28"
$ws.Cells.Item(21,2).Value = "14:41:11"
$ws.Cells.Item(21,3).Value = "14:00"
$ws.Cells.Item(21,4).Value = "Bathroom"
$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value = "88.7%"
$ws.Cells.Item(21,6).Value = "Active"
$ws.Cells.Item(22,1).NumberFormat = "@"
$ws.Cells.Item(22,1).Value = "2026-01-28"
$ws.Cells.Item(22,2).Value = "14:41:23"
$ws.Cells.Item(22,3).Value = "14:00"
$ws.Cells.Item(22,4).Value = "Bathroom"
$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value = "88.7%"
$ws.Cells.Item(22,6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-01-28"
$ws.Cells.Item(14,2).Value = "14:40:35"
$ws.Cells.Item(14,3).Value = "14:00"
$ws.Cells.Item(14,4).Value = "Bathroom"
$ws.Cells.Item(14,5).Value = "22.7C"
$ws.Cells.Item(14,6).Value = "Active"
$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-01-28"
$ws.Cells.Item(15,2).Value = "14:40:39"
$ws.Cells.Item(15,3).Value = "14:00"
$ws.Cells.Item(15,4).Value = "Bathroom"
$ws.Cells.Item(15,5).Value = "22.6C"
$ws.Cells.Item(15,6).Value = "Active"
$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-01-28"
$ws.Cells.Item(16,2).Value = "14:40:43"
$ws.Cells.Item(16,3).Value = "14:00"
$ws.Cells.Item(16,4).Value = "Bathroom"
$ws.Cells.Item(16,5).Value = "22.7C"
$ws.Cells.Item(16,6).Value = "Active"
$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = "2026-01-28"
$ws.Cells.Item(17,2).Value = "14:40:52"
$ws.Cells.Item(17,3).Value = "14:00"
$ws.Cells.Item(17,4).Value = "Bathroom"
$ws.Cells.Item(17,5).Value = "22.7C"
$ws.Cells.Item(17,6).Value = "Active"
$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = "2026-01-28"
$ws.Cells.Item(18,2).Value = "14:40:55"
$ws.Cells.Item(18,3).Value = "14:00"
$ws.Cells.Item(18,4).Value = "Bathroom"
$ws.Cells.Item(18,5).Value = "22.7C"
$ws.Cells.Item(18,6).Value = "Active"
$ws.Cells.Item(19,1).NumberFormat = "@"
$ws.Cells.Item(19,1).Value = "2026-01-28"
$ws.Cells.Item(19,2).Value = "14:41:00"
$ws.Cells.Item(19,3).Value = "14:00"
$ws.Cells.Item(19,4).Value = "Bathroom"
$ws.Cells.Item(19,5).Value = "22.7C"
$ws.Cells.Item(19,6).Value = "Active"
$ws.Cells.Item(20,1).NumberFormat = "@"
$ws.Cells.Item(20,1).Value = "2026-01-28"
$ws.Cells.Item(20,2).Value = "14:41:03"
$ws.Cells.Item(20,3).Value = "14:00"
$ws.Cells.Item(20,4).Value = "Bathroom"
$ws.Cells.Item(20,5).Value = "22.7C"
$ws.Cells.Item(20,6).Value = "Active"
$ws.Cells.Item(21,1).NumberFormat = "@"
$ws.Cells.Item(21,1).Value = "2026-01-28"
$ws.Cells.Item(21,2).Value = "14:41:11"
$ws.Cells.Item(21,3).Value = "14:00"
$ws.Cells.Item(21,4).Value = "Bathroom"
$ws.Cells.Item(21,5).Value = "22.7C"
$ws.Cells.Item(21,6).Value = "Active"
$ws.Cells.Item(22,1).NumberFormat = "@"
$ws.Cells.Item(22,1).Value = "2026-01-28"
$ws.Cells.Item(22,2).Value = "14:41:24"
$ws.Cells.Item(22,3).Value = "14:00"
$ws.Cells.Item(22,4).Value = "Bathroom"
$ws.Cells.Item(22,5).Value = "22.7C"
$ws.Cells.Item(22,6).Value = "Active"

$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = "2026-01-28"
$ws.Cells.Item(11,2).Value = "14:40:35"
$ws.Cells.Item(11,3).Value = "14:00"
$ws.Cells.Item(11,4).Value = "Living Room Main Door"
$ws.Cells.Item(11,5).Value = "ENTER"
$ws.Cells.Item(11,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = "2026-01-28"
$ws.Cells.Item(12,2).Value = "14:40:35"
$ws.Cells.Item(12,3).Value = "14:00"
$ws.Cells.Item(12,4).Value = "Bathroom Door"
$ws.Cells.Item(12,5).Value = "ENTER"
$ws.Cells.Item(12,6).Value = "User ENTERED Bathroom"
$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = "2026-01-28"
$ws.Cells.Item(13,2).Value = "14:40:38"
$ws.Cells.Item(13,3).Value = "14:00"
$ws.Cells.Item(13,4).Value = "Living Room Main Door"
$ws.Cells.Item(13,5).Value = "EXIT"
$ws.Cells.Item(13,6).Value = "User EXITED Living Room Main Door"
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = "2026-01-28"
$ws.Cells.Item(14,2).Value = "14:40:40"
$ws.Cells.Item(14,3).Value = "14:00"
$ws.Cells.Item(14,4).Value = "Bathroom Door"
$ws.Cells.Item(14,5).Value = "EXIT"
$ws.Cells.Item(14,6).Value = "User EXITED Bathroom"
$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = "2026-01-28"
$ws.Cells.Item(15,2).Value = "14:40:46"
$ws.Cells.Item(15,3).Value = "14:00"
$ws.Cells.Item(15,4).Value = "Living Room Main Door"
$ws.Cells.Item(15,5).Value = "ENTER"
$ws.Cells.Item(15,6).Value = "User ENTERED Living Room Main Door"
$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = "2026-01-28"
$ws.Cells.Item(16,2).Value = "14:40:46"
$ws.Cells.Item(16,3).Value = "14:00"
$ws.Cells.Item(16,4).Value = "Bathroom Door"
$ws.Cells.Item(16,5).Value = "ENTER"
$ws.Cells.Item(16,6).Value = "User ENTERED Bathroom"

$ws = $wb.Worksheets.Item("Camera")
$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "2026-01-28"
$ws.Cells.Item(5,2).Value = "14:40:35"
$ws.Cells.Item(5,3).Value = "14:00"
$ws.Cells.Item(5,4).Value = "Living Room Main Door"
$ws.Cells.Item(5,5).Value = "Image Captured"
$ws.Cells.Item(5,6).Value = "Active"
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "2026-01-28"
$ws.Cells.Item(6,2).Value = "14:40:37"
$ws.Cells.Item(6,3).Value = "14:00"
$ws.Cells.Item(6,4).Value = "Living Room Main Door"
$ws.Cells.Item(6,5).Value = "Image Captured"
$ws.Cells.Item(6,6).Value = "Active"
$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = "2026-01-28"
$ws.Cells.Item(7,2).Value = "14:40:40"
$ws.Cells.Item(7,3).Value = "14:00"
$ws.Cells.Item(7,4).Value = "Living Room Main Door"
$ws.Cells.Item(7,5).Value = "Image Captured"
$ws.Cells.Item(7,6).Value = "Active"
$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = "2026-01-28"
$ws.Cells.Item(8,2).Value = "14:40:47"
$ws.Cells.Item(8,3).Value = "14:00"
$ws.Cells.Item(8,4).Value = "Living Room Main Door"
$ws.Cells.Item(8,5).Value = "Image Captured"
$ws.Cells.Item(8,6).Value = "Active"
